$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture existing Sample_ID values (column B, rows 2-10) before we overwrite the header/column.
$lastRow = 10
$sampleIds = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $sampleIds[$r] = $ws.Cells.Item($r, 2).Text
}

# New header row: Number | Site | Zone | Date | Time | Replicate
$ws.Range("B1").Value = "Site"
$ws.Range("C1").Value = "Zone"
$ws.Range("D1").Value = "Date"
$ws.Range("E1").Value = "Time"
$ws.Range("F1").Value = "Replicate"

# Split each Sample_ID (e.g. TEMPEST_AqWell_20250904_1300_B) into its components.
for ($r = 2; $r -le $lastRow; $r++) {
    $parts = $sampleIds[$r].Split("_")
    $ws.Cells.Item($r, 2).Value = $parts[0]
    $ws.Cells.Item($r, 3).Value = $parts[1]
    $ws.Cells.Item($r, 4).Value = $parts[2]
    $ws.Cells.Item($r, 5).Value = $parts[3]
    $ws.Cells.Item($r, 6).Value = $parts[4]
}

$ws.Range("E1").Select()
